$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new animation rows (Talk / Laugh) following the existing table layout
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Talk"
$ws.Range("E7").Value = 183
$ws.Range("F7").Value = 328

$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "Laugh"
$ws.Range("E8").Value = 329
$ws.Range("F8").Value = 471

# Update selection to match the recorded UI state after the edit
$ws.Range("E13").Select()
